$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Give the matrix a title in the previously-empty corner cell.
$ws.Range("A1").Value = "Надёжность"

# 2. The title cell keeps the font/border used throughout the table, but
#    should NOT have the centered/wrapped alignment that the rest of the
#    table uses (it reverts to the default General/Bottom, no wrap).
$ws.Range("A1").WrapText = $false
$ws.Range("A1").VerticalAlignment = -4107   # xlBottom (default)

# 3. The numeric comparison matrix (B2:G7) and the priority-vector column
#    (H2:H7) become horizontally centered, in addition to the existing
#    vertical centering + wrap. (Applied to each area separately since
#    multi-area unions are not reliably supported.)
$ws.Range("B2:G7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2:G7").VerticalAlignment = -4108     # xlCenter
$ws.Range("B2:G7").WrapText = $true

$ws.Range("H2:H7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H2:H7").VerticalAlignment = -4108     # xlCenter
$ws.Range("H2:H7").WrapText = $true

# 4. The empty bordered box (A8:G10, merged) loses the centered/wrapped
#    alignment as well, matching the rest of the "plain" bordered cells.
#    Copy the (now plain) formatting from A1 in a single paste so the
#    whole merged area converges directly onto that same style.
$ws.Range("A1").Copy()
$ws.Range("A8:G10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Headers (B1:H1), row labels (A2:A7) and the summary text cells (H8:H10)
# keep their existing vertical-centered / wrapped formatting untouched.
